$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F (shifts old F->G, G->H, H->I, I->J and adjusts formulas)
$ws.Columns("F:F").Insert()

# The inserted column inherits formatting from column E; strip it so the new
# column F cells carry no style, matching a freshly pasted-in data column.
$ws.Range("F3:F118").ClearFormats()

# Row 1 / Row 2 headers: column F is a brand-new "FIOC" data column, and the UKDataID
# row (row 1) across the whole data block is now uniformly "BoE data".
$ws.Range("B1:J1").Value = "BoE data"
$ws.Range("F2").Value = "FIOC"

# New FIOC data values for rows 3:118 (years 1900-2015)
$fioc = @(0.6480795886216052,0.70611139818092794,0.69362147672977004,0.70897589941903383,0.7174744476015551,0.64844158551087128,0.58972586199271648,0.55851855848571053,0.66495798366551995,0.64738900790244003,0.61681297339865837,0.60248414620920743,0.59420695044603289,0.57369373598414697,0.66925828634811135,1.1114328485162206,1.0465471038380851,0.98945748605150163,0.98531654935776836,0.85525952691875951,0.75769471386012388,0.80080408214870935,0.78105283318857166,0.80159262647952312,0.7889262761317829,0.73936013831794878,0.80434990428522901,0.72452885558871827,0.73726582541027685,0.73029991013064299,0.73611605623927989,0.85310912285460272,0.8574858979257034,0.83721003260747828,0.78113964044109507,0.7720537240922688,0.76650515029525212,0.75054072438921571,0.77186179238872088,0.86290305617071417,1.0899634088997714,1.104772609246258,1.1147570913896014,1.0913779185254497,1.1708494569147663,1.1943352888034577,0.88399572438019269,0.7776221744606997,0.7339289973406925,0.73928649625009479,0.7261664479403106,0.69912014297784109,0.70212258280882334,0.69631228190033545,0.7031340312570008,0.68410036412747677,0.63989134840182038,0.63189172166214358,0.62720447150599112,0.6406608348737548,0.60886789268656061,0.59252280976465843,0.60956069667838009,0.58604035759740469,0.55854298436134153,0.55328258794840246,0.54956055499138401,0.55715617104891213,0.52783335053790725,0.48355734902907782,0.41453605514388675,0.41915758009558274,0.44112097409926526,0.4084360420131109,0.48969016021640993,0.5388470766123884,0.50361904308932792,0.46805992005065922,0.44303014041719857,0.43879391420144609,0.47605283753390981,0.47987316463702606,0.47157294518147608,0.4484033364313198,0.43345569891561891,0.41690274986603559,0.43154175336832234,0.39381876518473902,0.37204910267948893,0.33915996819860178,0.32709617327019769,0.32460079564031213,0.34858257894374162,0.34225891608512676,0.32231346291016832,0.33757078740972812,0.3229649104991788,0.32307904628497119,0.27336557725008404,0.34155887558663817,0.3192656059738615,0.29999517561767064,0.28519393441207214,0.27003476418544575,0.2622845237852483,0.24091926540034331,0.26231770228847717,0.24104185834990891,0.35740023014945832,0.39317497619204939,0.33925223864010529,0.30296463797257989,0.37376197317470028,0.39445433661578017,0.35925101974799106,0.32083202653607018)
for ($i = 0; $i -lt $fioc.Length; $i++) {
    $ws.Cells.Item(3 + $i, 6).Value = $fioc[$i]
}

# Restore the selection to match the edited column
$ws.Range("F1:F1048576").Select()
